# Attendance sheet update: add "May 30" column of attendance checkmarks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# Add the day number 30 to the new date column (K7), continuing 25,26,27,28,29.
$ws.Range("K7").Value = 30

# Mark attendance ("checked-in") for day 30 for the people who attended:
# row 8 (Abdullah), 13 (Rabiha), 14 (Sanwal), 16 (Areej), 17 (Husna), 18 (Tayeba), 19 (Vusqa)
$ws.Range("K8").Value = "✓"
$ws.Range("K13").Value = "✓"
$ws.Range("K14").Value = "✓"
$ws.Range("K16").Value = "✓"
$ws.Range("K17").Value = "✓"
$ws.Range("K18").Value = "✓"
$ws.Range("K19").Value = "✓"

# Re-apply the per-row "Total presence" COUNTA formula across Q8:Q21 so Excel
# folds the now-identical formulas back into a single shared-formula group
# (matches the template behaviour of filling the formula down the column).
$ws.Range("Q8:Q21").Formula = "=COUNTA(F8:P8)"

$wb.Save()
